$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.188.61"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "2.260.84"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'113.17"
$ws.Range("E5").Value = "  +4.74%  "
$ws.Range("D6").Value = "'264.77"
$ws.Range("E6").Value = "  -2.37%  "
$ws.Range("D7").Value = "'0.617"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.599"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").Value = "'47.69"
$ws.Range("E10").Value = "  +1.31%  "
$ws.Range("D11").Value = "'0.0923"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "'8.74"
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("D14").Value = "'15.41"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").Value = "2.601.79"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "'0.852"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "2.264.08"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "43.132.36"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "'6.97"
$ws.Range("E20").Value = "  +10.60%  "
$ws.Range("D21").Value = "'70.90"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").Value = "'2.39"
$ws.Range("E22").Value = "  -3.85%  "
$ws.Range("D23").Value = "'9.73"
$ws.Range("E23").Value = "  +4.42%  "
$ws.Range("D24").Value = "'229.84"
$ws.Range("E24").Value = "  -1.64%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'11.28"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'3.88"
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").Value = "'41.19"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "'3.39"
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "'171.44"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").Value = "'21.18"
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("D34").Value = "'0.0897"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").Value = "'5.54"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("D37").Value = "'4.61"
$ws.Range("E37").Value = "  -5.93%  "
$ws.Range("D38").Value = "'0.0349"
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +2.98%  "
$ws.Range("E40").Value = "  -8.67%  "
$ws.Range("D41").Value = "'14.19"
$ws.Range("E41").Value = "  +16.41%  "
$ws.Range("D42").Value = "'74.93"
$ws.Range("E42").Value = "  +12.92%  "
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").Value = "'0.233"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "'6.12"
$ws.Range("E45").Value = "  +11.04%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'1.36"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("D49").Value = "'0.0985"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'99.70"
$ws.Range("E51").Value = "  +0.35%  "
